# Update a couple of values in the OrangeHRM_Excel workbook.
$wb = $excel.ActiveWorkbook

# --- Admin sheet: Username column (D2) ---
$adminSheet = $wb.Worksheets.Item("Admin")
$adminSheet.Range("D2").Value = "4482716A"

# --- Jira sheet: Error description column (B2) ---
$jiraSheet = $wb.Worksheets.Item("Jira")
$jiraSheet.Range("B2").Value = "5FA3C2312892FD51F30690CB47131C4C"
